# "Commiting new test cases"
# TestCase_F1 (row 2) and TestCase_F2 (row 3) results change from SKIP to FAIL,
# TestCase_F3 (row 4) result changes from SKIP to PASS, on the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "FAIL"
$ws.Range("D3").Value = "FAIL"
$ws.Range("D4").Value = "PASS"
